$wb = $excel.ActiveWorkbook

# Sheet names that need updating: 展览 (Exhibition) and 全部类型 (All types)
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1791
    $ws.Range("F3").Value = 8153
    $ws.Range("F5").Value = 297
}
